# Generate Report for Handoff
# Replaces the "in sync with en-US / handed back" report data with a fresh
# "ready for handoff" report (new source GUID file names, new status,
# new timestamps) and drops the now-redundant "Latest Target File" /
# "Latest Handback File" columns (F, G) from the per-language detail sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# New values used throughout the workbook
# ---------------------------------------------------------------------
$newMdA      = "59a5d24b-0e48-4de4-987d-ad457091dcba.md"
$newMdB      = "ffff66f9815e-fd99-44a1-aa5a-a27856595390.md"
$newStatus   = "Ready for handoff"
$newOverviewDate = "2016-48-17 12:48:10"

$newZhXlf    = "59a5d24b-0e48-4de4-987d-ad457091dcba.4adfb3cae53adff77fdfb5fa97466f6bbe59c6c6.zh-cn.xlf"
$newZhHandoffDt = "2016-03-17 12:48:06"
$newZhHandbackDt = "0001-01-01 00:00:00"

$newDeXlf    = "59a5d24b-0e48-4de4-987d-ad457091dcba.4adfb3cae53adff77fdfb5fa97466f6bbe59c6c6.de-de.xlf"
$newDeHandoffDt = "2016-03-17 12:48:10"
$newDeHandbackDt = "0001-01-01 00:00:00"

$ghRootA = "https://github.com/OpenLocalizationTest/oltest/blob/fb1e20dff7e0fb356b4f05c7d2c7a0db7077d9b3/e2e/$newMdA"
$ghRootB = "https://github.com/OpenLocalizationTest/oltest/blob/fb1e20dff7e0fb356b4f05c7d2c7a0db7077d9b3/e2e/$newMdB"

$ghZhXlfA = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f9759d60283fd25109cfa8c45e6cae76c2af1626/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newZhXlf"
$ghDeXlfA = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dc9a837fb6249d3a49be1fd58f4bed710ebbb61a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newDeXlf"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(2,1).Value2 = $newMdA
$wsOverview.Cells.Item(2,2).Value2 = $newStatus
$wsOverview.Cells.Item(2,3).Value2 = $newStatus
$wsOverview.Cells.Item(2,4).Value2 = $newOverviewDate

$wsOverview.Cells.Item(3,1).Value2 = $newMdB
$wsOverview.Cells.Item(3,2).Value2 = $newStatus
$wsOverview.Cells.Item(3,3).Value2 = $newStatus
$wsOverview.Cells.Item(3,4).Value2 = $newOverviewDate

# Rebuild hyperlinks on the Overview sheet (A2 / A3 point at the .md files)
$wsOverview.Range("A1").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item(2,1), $ghRootA, "", "", $newMdA)
$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item(3,1), $ghRootB, "", "", $newMdB)

# ---------------------------------------------------------------------
# Helper: update one language detail sheet (zh-cn / de-de)
# ---------------------------------------------------------------------
function Update-LangSheet($ws, $xlf, $handoffDt, $handbackDt, $ghXlf) {
    $ws.Cells.Item(2,1).Value2 = $newMdA
    $ws.Cells.Item(2,2).Value2 = ".md"
    $ws.Cells.Item(2,3).Value2 = $newStatus
    $ws.Cells.Item(2,4).Value2 = $xlf
    $ws.Cells.Item(2,5).Value2 = $handoffDt
    $ws.Cells.Item(2,8).Value2 = $handbackDt
    $ws.Cells.Item(2,9).Value2 = "Include"

    $ws.Cells.Item(3,1).Value2 = $newMdB
    $ws.Cells.Item(3,2).Value2 = ".md"
    $ws.Cells.Item(3,3).Value2 = $newStatus
    $ws.Cells.Item(3,4).Value2 = $xlf
    $ws.Cells.Item(3,5).Value2 = $handoffDt
    $ws.Cells.Item(3,8).Value2 = $handbackDt
    $ws.Cells.Item(3,9).Value2 = "Include"

    # Drop the "Latest Target File" (F) and "Latest Handback File" (G) columns
    # for rows 2 and 3 - they are no longer populated in this report.
    $ws.Range("F2:G3").ClearContents()
    $ws.Range("F2:G3").ClearFormats()

    # Rebuild the hyperlinks collection: A/B/D for row 2 and row 3 only.
    $ws.Range("A1").Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Cells.Item(2,1), $ghRootA, "", "", $newMdA)
    $ws.Hyperlinks.Add($ws.Cells.Item(2,2), $ghRootA, "", "", ".md")
    $ws.Hyperlinks.Add($ws.Cells.Item(2,4), $ghXlf, "", "", $xlf)
    $ws.Hyperlinks.Add($ws.Cells.Item(3,1), $ghRootB, "", "", $newMdB)
    $ws.Hyperlinks.Add($ws.Cells.Item(3,2), $ghRootB, "", "", ".md")
    $ws.Hyperlinks.Add($ws.Cells.Item(3,4), $ghXlf, "", "", $xlf)
}

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
Update-LangSheet $wsZh $newZhXlf $newZhHandoffDt $newZhHandbackDt $ghZhXlfA

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
Update-LangSheet $wsDe $newDeXlf $newDeHandoffDt $newDeHandbackDt $ghDeXlfA
